$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.797.89'
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").Value = '2.355.22'

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.691'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +5.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.27'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '77.19'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +4.85%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.626'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +20.11%  '

$ws.Range("E10").Value = '  +3.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.36'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '33.86'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +21.97%  '

$ws.Range("E13").Value = '  +19.25%  '

$ws.Range("E14").Value = '  +1.84%  '

$ws.Range("D15").Value = '2.704.27'
$ws.Range("E15").Value = '  -0.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.94'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.926'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.97%  '

$ws.Range("D18").Value = '2.352.84'
$ws.Range("E18").Value = '  -0.78%  '

$ws.Range("D19").Value = '43.763.15'
$ws.Range("E19").Value = '  +0.43%  '

$ws.Range("E20").Value = '  +2.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.65'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.53'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.12'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.54%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("E25").Value = '  +2.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +7.87%  '

$ws.Range("E27").Value = '  -4.90%  '

$ws.Range("E28").Value = '  +16.38%  '

$ws.Range("E29").Value = '  +2.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.04'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.91'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.17%  '

$ws.Range("E32").Value = '  -4.30%  '

$ws.Range("E33").Value = '  +5.86%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0757'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +8.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.79%  '

$ws.Range("E36").Value = '  +6.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.79'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("E38").Value = '  -1.13%  '

$ws.Range("E39").Value = '  -3.32%  '

$ws.Range("E40").Value = '  +7.32%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '19.49'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("E42").Value = '  +1.16%  '

$ws.Range("E43").Value = '  +15.10%  '

$ws.Range("E45").Value = '  +7.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.53'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +11.80%  '

$ws.Range("E47").Value = '  +3.35%  '

$ws.Range("E48").Value = '  +1.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.52'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.93%  '

$ws.Range("E50").Value = '  -0.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.51'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.48%  '
